$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "is_locked_lbl" (D1) and "is_enabled_lbl" (E1) template columns,
# shifting the remaining columns (order_by, rem) left so they become D1/E1,
# leaving the row with only columns A1:E1 (F1/G1 no longer exist).
$ws.Range("D1:E1").Delete(-4159) | Out-Null
